{"js": "// The document contains a single 20-row x 5-column table of two-digit\n// \u00f7 one-digit division prompts (\"NN\u00f7N=\"). Five of the rows (0, 4, 8,\n// 12, 16) hold the actual problems; the rows in between are left blank\n// for answers. This script overwrites the text of each of those 25\n// populated cells in place (row-major order), leaving every other part\n// of the document (including per-run formatting) untouched.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// row index -> new cell text values (left to right)\nconst newRowValues = {\n  0: [\"82\u00f75=\", \"94\u00f73=\", \"80\u00f77=\", \"51\u00f74=\", \"67\u00f79=\"],\n  4: [\"29\u00f76=\", \"39\u00f78=\", \"53\u00f75=\", \"82\u00f75=\", \"61\u00f77=\"],\n  8: [\"79\u00f79=\", \"19\u00f72=\", \"54\u00f74=\", \"76\u00f73=\", \"57\u00f74=\"],\n  12: [\"89\u00f79=\", \"87\u00f73=\", \"82\u00f78=\", \"45\u00f78=\", \"81\u00f75=\"],\n  16: [\"61\u00f73=\", \"78\u00f73=\", \"64\u00f75=\", \"25\u00f74=\", \"37\u00f75=\"],\n};\n\nfor (const rowIndexStr of Object.keys(newRowValues)) {\n  const rowIndex = Number(rowIndexStr);\n  const values = newRowValues[rowIndex];\n  for (let col = 0; col < values.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    const range = cell.body.paragraphs.getFirst().getRange();\n    range.insertText(values[col], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single 20-row x 5-column table of two-digit\n# / one-digit division prompts (\"NN/N=\"). Five of the rows (the 1st,\n# 5th, 9th, 13th, and 17th, 1-based) hold the actual problems; the rows\n# in between are left blank for answers. This script overwrites the\n# text of each of those 25 populated cells in place (row by row, left\n# to right), leaving every other part of the document (including\n# per-run formatting) untouched.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# 1-based table row -> new cell text values (left to right)\n$newRowValues = @{\n    1  = @(\"82\u00f75=\", \"94\u00f73=\", \"80\u00f77=\", \"51\u00f74=\", \"67\u00f79=\")\n    5  = @(\"29\u00f76=\", \"39\u00f78=\", \"53\u00f75=\", \"82\u00f75=\", \"61\u00f77=\")\n    9  = @(\"79\u00f79=\", \"19\u00f72=\", \"54\u00f74=\", \"76\u00f73=\", \"57\u00f74=\")\n    13 = @(\"89\u00f79=\", \"87\u00f73=\", \"82\u00f78=\", \"45\u00f78=\", \"81\u00f75=\")\n    17 = @(\"61\u00f73=\", \"78\u00f73=\", \"64\u00f75=\", \"25\u00f74=\", \"37\u00f75=\")\n}\n\nforeach ($rowIndex in $newRowValues.Keys) {\n    $values = $newRowValues[$rowIndex]\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $cell = $table.Cell($rowIndex, $col)\n        $cell.Range.Text = $values[$col - 1]\n    }\n}\n"}
